# Updates the "cryptos" price/volume table to reflect the latest scrape.
# D (Price) values are written with a leading apostrophe so Excel keeps
# them as plain text (matching the inline strings already in the sheet)
# instead of re-interpreting dotted numbers such as "3.750" or
# "29.081.58" as numeric/date values and mangling their formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.081.58'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  -1.28%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.971.30'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  -0.91%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '1.015'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  +0.81%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '329.23'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -0.22%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '1.013'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +0.66%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4964'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -0.29%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.4196'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -0.32%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '54.35'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +4.50%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.09317'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +4.83%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '1.097'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -2.00%  '; ForceText = $false }
    @{ Cell = 'B12'; Value = 'Solana'; ForceText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; ForceText = $false }
    @{ Cell = 'D12'; Value = '22.75'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'B13'; Value = 'WrappedEther'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $false }
    @{ Cell = 'D13'; Value = '2.014.13'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +1.46%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '7.876'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -2.26%  '; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -0.72%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '1.016'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +0.93%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '0.00001110'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +0.50%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '91.66'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -4.61%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.06704'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +1.32%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '19.14'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -2.89%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '1.012'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +0.71%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '5.945'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '29.113.10'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -1.21%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '11.94'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.76%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '2.267'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.26%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '2.217.00'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '20.74'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +0.88%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '156.79'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.39%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '6.227'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -4.76%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '2.264'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.97%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '127.34'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -0.39%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  -0.79%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.09824'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -1.16%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.497'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -4.49%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '5.802'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.58%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '3.750'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -1.10%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.02413'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -1.55%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '1.325'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +3.03%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.06397'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.78%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '9.022'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -6.06%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '0.6465'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -0.73%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -2.17%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.1998'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -3.28%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '1.012'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.6187'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -2.55%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '1.355'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +6.80%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '13.33'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -0.34%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.174'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '3.493'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -1.10%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.00000000332'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.06962'; ForceText = $true }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).Value = "'" + $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

# Strip the quote-prefix style the apostrophes above introduced so the
# Price column cells keep their original (unstyled / General) look.
$ws.Range("D2:D51").ClearFormats()
